$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" (E), "Valor Mora" (F) and "Salario Basico" (G) columns
# for the EC (Estado de Cuenta) detail rows 16-22 with the refreshed database values.
# Periods are re-sorted ascending (2103 .. 2201) and Salario Basico is updated
# uniformly to the new base salary value for every period row.

$ws.Range("E16").Value = "2103"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 877803

$ws.Range("E17").Value = "2104"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 877803

$ws.Range("E18").Value = "2109"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 877803

$ws.Range("E19").Value = "2110"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = 877803

$ws.Range("E20").Value = "2111"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 877803

$ws.Range("E21").Value = "2112"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 877803

$ws.Range("E22").Value = "2201"
$ws.Range("F22").Value = 29260
$ws.Range("G22").Value = 877803
